# Update "F" column (想去人数 / "want to go" headcount) values on both the
# "展览" and "全部类型" sheets, which mirror each other's data.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 42
    8  = 71
    10 = 1179
    11 = 1483
    12 = 330
    13 = 368
    15 = 126
    19 = 263
    20 = 285
    21 = 312
    22 = 1694
    25 = 170
    26 = 640
    28 = 230
    29 = 4053
    31 = 474
    32 = 248
    33 = 1046
    34 = 124
    36 = 316
    38 = 165
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
